$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($rangeAddr, $value) {
    $c = $ws.Range($rangeAddr)
    $c.NumberFormat = "@"
    $c.Value = $value
    $c.Style = "Normal"
}

Set-TextValue "D2" "317.05"
Set-TextValue "E2" "2.28%"

Set-TextValue "D3" "41.09"
Set-TextValue "E3" "-0.05%"

Set-TextValue "D4" "5.142"
Set-TextValue "E4" "0.40%"

Set-TextValue "D5" "0.07636"
Set-TextValue "E5" "-0.65%"

$ws.Range("B6").Value = "FTXToken"
$ws.Range("C6").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
Set-TextValue "D6" "1.680"
Set-TextValue "E6" "3.39%"

$ws.Range("B7").Value = "MXToken"
$ws.Range("C7").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
Set-TextValue "D7" "0.9337"
Set-TextValue "E7" "1.35%"

$ws.Range("B8").Value = "BTSEToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
Set-TextValue "D8" "2.425"
Set-TextValue "E8" "-1.74%"

$ws.Range("B9").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C9").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
Set-TextValue "D9" "0.1245"
Set-TextValue "E9" "1.95%"

$ws.Range("B10").Value = "WazirX"
$ws.Range("C10").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
Set-TextValue "D10" "0.1826"
Set-TextValue "E10" "0.17%"

$ws.Range("B11").Value = "MandalaExchangeToken"
$ws.Range("C11").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
Set-TextValue "D11" "0.09083"
Set-TextValue "E11" "-0.85%"

$ws.Range("B12").Value = "BitrueCoin"
$ws.Range("C12").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
Set-TextValue "D12" "0.04142"
Set-TextValue "E12" "-4.37%"

$ws.Range("B13").Value = "BitMartToken"
$ws.Range("C13").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
Set-TextValue "D13" "0.1056"
Set-TextValue "E13" "0.48%"

$ws.Range("B14").Value = "BitForexToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
Set-TextValue "D14" "0.001291"
Set-TextValue "E14" "4.77%"

$ws.Range("B15").Value = "TigerCash"
$ws.Range("C15").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
Set-TextValue "D15" "0.005933"
Set-TextValue "E15" "2.23%"

$ws.Range("B16").Value = "UpBots"
$ws.Range("C16").Value = "https://coinranking.com/coin/m5ozaAIK6+upbots-ubxt"
Set-TextValue "D16" "0.007491"
Set-TextValue "E16" "1,897.31%"

$ws.Range("B17").Value = "LEO"
$ws.Range("C17").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
Set-TextValue "D17" "3.351"
Set-TextValue "E17" "-0.07%"

$ws.Range("B18").Value = "GateToken"
$ws.Range("C18").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
Set-TextValue "D18" "4.320"
Set-TextValue "E18" "0.81%"

Set-TextValue "E19" "1.46%"

Set-TextValue "D20" "8.408"
Set-TextValue "E20" "21.88%"

Set-TextValue "E21" "-2.83%"

Set-TextValue "D22" "0.2869"
Set-TextValue "E22" "7.29%"

Set-TextValue "D23" "0.04043"
Set-TextValue "E23" "0.19%"

Set-TextValue "D24" "0.001273"
Set-TextValue "E24" "0.80%"

Set-TextValue "D25" "0.004084"
Set-TextValue "E25" "0.16%"

Set-TextValue "E26" "0.46%"

Set-TextValue "D38" "0.02504"
Set-TextValue "E38" "1.87%"

Set-TextValue "D39" "0.05247"
Set-TextValue "E39" "-0.24%"

Set-TextValue "D40" "0.007785"
Set-TextValue "E40" "-0.63%"

Set-TextValue "D41" "0.1298"
Set-TextValue "E41" "-1.22%"

Set-TextValue "D42" "0.007074"
Set-TextValue "E42" "4.13%"

Set-TextValue "D43" "0.002078"
Set-TextValue "E43" "12.84%"

Set-TextValue "D44" "0.008241"

Set-TextValue "E45" "2.30%"

Set-TextValue "D46" "0.00006685"
Set-TextValue "E46" "-1.45%"

Set-TextValue "E47" "0.47%"

Set-TextValue "D48" "0.2233"
Set-TextValue "E48" "-0.80%"

Set-TextValue "D49" "0.004216"
Set-TextValue "E49" "2.95%"

Set-TextValue "D50" "0.00002108"
Set-TextValue "E50" "0.47%"

Set-TextValue "D51" "0.0002007"
Set-TextValue "E51" "0.47%"
